# Exercice4: add "langue" column (E) to Sheet1, shift densite to F,
# and add a new "population" worksheet summarising language stats.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Preserve the existing "densite" column (values + header style) by
# copying E1:E54 -> F1:F54 before E is overwritten with language data.
$ws1.Range("E1:E54").Copy($ws1.Range("F1"))

# New header for column E.
$ws1.Range("E1").Value = "langue"

# Per-country language map (row 2 = Italy ... row 54 = Cyprus).
$langues = @{
    2 = '{''ita'': ''Italian'', ''cat'': ''Catalan''}'
    3 = '{''eng'': ''English'', ''fra'': ''French'', ''nrf'': ''Jèrriais''}'
    4 = '{''mkd'': ''Macedonian''}'
    5 = '{''lav'': ''Latvian''}'
    6 = '{''est'': ''Estonian''}'
    7 = '{''bel'': ''Belarusian'', ''rus'': ''Russian''}'
    8 = '{''fra'': ''French'', ''gsw'': ''Swiss German'', ''ita'': ''Italian'', ''roh'': ''Romansh''}'
    9 = '{''sqi'': ''Albanian'', ''srp'': ''Serbian''}'
    10 = '{''deu'': ''German''}'
    11 = '{''deu'': ''German'', ''fra'': ''French'', ''nld'': ''Dutch''}'
    12 = '{''isl'': ''Icelandic''}'
    13 = '{''slv'': ''Slovene''}'
    14 = '{''swe'': ''Swedish''}'
    15 = '{''ukr'': ''Ukrainian''}'
    16 = '{''eng'': ''English'', ''fra'': ''French'', ''nfr'': ''Guernésiais''}'
    17 = '{''ces'': ''Czech'', ''slk'': ''Slovak''}'
    18 = '{''ita'': ''Italian'', ''lat'': ''Latin''}'
    19 = '{''cnr'': ''Montenegrin''}'
    20 = '{''deu'': ''German''}'
    21 = '{''pol'': ''Polish''}'
    22 = '{''fra'': ''French''}'
    23 = '{''slk'': ''Slovak''}'
    24 = '{''lit'': ''Lithuanian''}'
    25 = '{''fra'': ''French''}'
    26 = '{''spa'': ''Spanish'', ''cat'': ''Catalan'', ''eus'': ''Basque'', ''glc'': ''Galician''}'
    27 = '{''nno'': ''Norwegian Nynorsk'', ''nob'': ''Norwegian Bokmål'', ''smi'': ''Sami''}'
    28 = '{''eng'': ''English'', ''gle'': ''Irish''}'
    29 = '{''fin'': ''Finnish'', ''swe'': ''Swedish''}'
    30 = '{''deu'': ''German''}'
    31 = '{''eng'': ''English''}'
    32 = '{''eng'': ''English''}'
    33 = '{''swe'': ''Swedish''}'
    34 = '{''dan'': ''Danish'', ''fao'': ''Faroese''}'
    35 = '{''eng'': ''English'', ''glv'': ''Manx''}'
    36 = '{''dan'': ''Danish''}'
    37 = '{''eng'': ''English'', ''mlt'': ''Maltese''}'
    38 = '{''ron'': ''Romanian''}'
    39 = '{''cat'': ''Catalan''}'
    40 = '{''ell'': ''Greek''}'
    41 = '{''hrv'': ''Croatian''}'
    42 = '{''ita'': ''Italian''}'
    43 = '{''nld'': ''Dutch''}'
    44 = '{''bul'': ''Bulgarian''}'
    45 = '{''nor'': ''Norwegian''}'
    46 = '{''deu'': ''German'', ''fra'': ''French'', ''ltz'': ''Luxembourgish''}'
    47 = '{''rus'': ''Russian''}'
    48 = '{''por'': ''Portuguese''}'
    49 = '{''hun'': ''Hungarian''}'
    50 = '{''ron'': ''Romanian''}'
    51 = '{''sqi'': ''Albanian''}'
    52 = '{''bos'': ''Bosnian'', ''hrv'': ''Croatian'', ''srp'': ''Serbian''}'
    53 = '{''srp'': ''Serbian''}'
    54 = '{''ell'': ''Greek'', ''tur'': ''Turkish''}'
}

foreach ($row in $langues.Keys) {
    $ws1.Range("E$row").Value = $langues[$row]
}

# New "population" worksheet, inserted right after Sheet1.
$wsPop = $wb.Worksheets.Add($null, $ws1)
$wsPop.Name = "population"

# Reuse the existing bold/bordered header style from Sheet1 instead of
# re-deriving an equivalent (but duplicate) style via direct formatting.
$ws1.Range("A1:C1").Copy($wsPop.Range("A1"))
$wsPop.Range("A1").Value = "langue"
$wsPop.Range("B1").Value = "nb_pays"
$wsPop.Range("C1").Value = "population_totale"

$wsPop.Range("A2").Value = "rus"
$wsPop.Range("B2").Value = 2
$wsPop.Range("C2").Value = 155137605

$wsPop.Range("A3").Value = "cat"
$wsPop.Range("B3").Value = 3
$wsPop.Range("C3").Value = 108331988

$wsPop.Range("A4").Value = "deu"
$wsPop.Range("B4").Value = 5
$wsPop.Range("C4").Value = 105240604

# Restore Sheet1 as the active tab (unchanged bookViews in the target).
$ws1.Activate()

